$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.143.13"
$ws.Range("E2").Value = "  +0.60%  "

# Row 3
$ws.Range("D3").Value = "1.749.74"
$ws.Range("E3").Value = "  -0.08%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").Value = "'236.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.58%  "

# Row 6
$ws.Range("D6").Value = "'1.0000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.10%  "

# Row 7
$ws.Range("D7").Value = "'0.5303"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.70%  "

# Row 8
$ws.Range("E8").Value = "  -1.61%  "

# Row 9
$ws.Range("D9").Value = "'0.06174"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.34%  "

# Row 10
$ws.Range("D10").Value = "1.746.17"
$ws.Range("E10").Value = "  -0.21%  "

# Row 11
$ws.Range("D11").Value = "'0.07175"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.99%  "

# Row 12
$ws.Range("D12").Value = "'15.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.61%  "

# Row 13
$ws.Range("D13").Value = "'0.6475"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.05%  "

# Row 14
$ws.Range("D14").Value = "'4.629"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.12%  "

# Row 15
$ws.Range("D15").Value = "'78.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.98%  "

# Row 16
$ws.Range("D16").Value = "'0.9998"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.08%  "

# Row 17
$ws.Range("E17").Value = "  +0.21%  "

# Row 18
$ws.Range("D18").Value = "26.042.43"
$ws.Range("E18").Value = "  +0.20%  "

# Row 19
$ws.Range("D19").Value = "'11.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.13%  "

# Row 20
$ws.Range("D20").Value = "'0.000006774"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.90%  "

# Row 21
$ws.Range("D21").Value = "1.971.80"
$ws.Range("E21").Value = "  -0.26%  "

# Row 22
$ws.Range("D22").Value = "'4.341"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.29%  "

# Row 23
$ws.Range("D23").Value = "'8.729"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.70%  "

# Row 24
$ws.Range("D24").Value = "'5.239"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.32%  "

# Row 25
$ws.Range("D25").Value = "'139.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.05%  "

# Row 26
$ws.Range("E26").Value = "  +1.31%  "

# Row 27
$ws.Range("D27").Value = "'15.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.96%  "

# Row 28
$ws.Range("D28").Value = "'1.802"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.66%  "

# Row 29
$ws.Range("D29").Value = "'104.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.50%  "

# Row 30
$ws.Range("D30").Value = "'0.08306"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.28%  "

# Row 31
$ws.Range("D31").Value = "'3.770"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.93%  "

# Row 32
$ws.Range("D32").Value = "'3.636"
$ws.Range("D32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'0.04637"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.76%  "

# Row 34
$ws.Range("D34").Value = "'2.642"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.19%  "

# Row 35
$ws.Range("D35").Value = "'1.007"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.74%  "

# Row 36
$ws.Range("D36").Value = "'0.6300"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.90%  "

# Row 37
$ws.Range("D37").Value = "'2.708"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.67%  "

# Row 38
$ws.Range("D38").Value = "'0.01621"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.61%  "

# Row 39
$ws.Range("D39").Value = "'1.971"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.68%  "

# Row 40
$ws.Range("D40").Value = "'0.9996"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.16%  "

# Row 41
$ws.Range("D41").Value = "'101.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.67%  "

# Row 42
$ws.Range("D42").Value = "'0.3927"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.08%  "

# Row 43
$ws.Range("D43").Value = "'0.7542"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.41%  "

# Row 44
$ws.Range("D44").Value = "'5.082"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.62%  "

# Row 45
$ws.Range("E45").Value = "  +2.82%  "

# Row 46
$ws.Range("D46").Value = "'6.343"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.35%  "

# Row 47
$ws.Range("D47").Value = "'0.05343"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.43%  "

# Row 48
$ws.Range("E48").Value = "  +2.96%  "

# Row 49
$ws.Range("D49").Value = "'30.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.90%  "

# Row 50
$ws.Range("E50").Value = "  +1.12%  "

# Row 51
$ws.Range("D51").Value = "'7.576"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.71%  "
